$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2023-07-03 15:50:05"
$ws.Range("F2").Value = "OUT"
